# Apply the compliance-results refresh (20.08 -> 20.16 total assets, and
# associated dependent metrics) across the relevant worksheets, plus two
# column-width tweaks on the "Expenses" columns of the diversification sheets.

$wb = $excel.ActiveWorkbook

# --- 12d1_Other_Investment_Companies: Total Assets (C2) ---
$ws = $wb.Worksheets.Item("12d1_Other_Investment_Companies")
$ws.Range("C2").Value = 20.16

# --- 12d2_Insurance_Companies: Total Assets (D2) ---
$ws = $wb.Worksheets.Item("12d2_Insurance_Companies")
$ws.Range("D2").Value = 20.16

# --- 12d3_Securities_Business: Total Assets (J2) ---
$ws = $wb.Worksheets.Item("12d3_Securities_Business")
$ws.Range("J2").Value = 20.16

# --- 40Act_Diversification ---
$ws = $wb.Worksheets.Item("40Act_Diversification")
# Narrow the "Expenses" column (K) from 21 to 10 characters wide.
$ws.Columns("K").ColumnWidth = 9.166666666666666
$ws.Range("H2").Value = 20.16
$ws.Range("J2").Value = 20.16
$ws.Range("K2").Value = 0.1512
$ws.Range("N2").Value = "AMGN, AMGN, VZ, VZ, PG, PG, PG, MRK, MRK, MRK, MCD, MCD, MCD, KO, KO, KO, JNJ, JNJ, JNJ, IBM, IBM, IBM, CVX, CVX, CVX, CSCO, CSCO, CSCO, AMGN, VZ"

# --- IRS_Diversification ---
$ws = $wb.Worksheets.Item("IRS_Diversification")
# Narrow the "expenses" column (J) from 21 to 10 characters wide.
$ws.Columns("J").ColumnWidth = 9.166666666666666
$ws.Range("H2").Value = 20.16
$ws.Range("J2").Value = 0.1512
$ws.Range("K2").Value = 1.008

# --- IRC_Diversification: Total Assets (K2) ---
$ws = $wb.Worksheets.Item("IRC_Diversification")
$ws.Range("K2").Value = 20.16

# --- Illiquid: Total Assets (C2) ---
$ws = $wb.Worksheets.Item("Illiquid")
$ws.Range("C2").Value = 20.16
